$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit per-cell updates mirroring the authoritative diff.
# Cells whose new text looks like a plain number (single decimal point) are
# written with the Text number format first so Excel keeps them as literal
# strings (matching the source inlineStr cells) instead of silently coercing
# them into numeric values; the style is reset back to Normal immediately
# after so no visible/persistent formatting change is introduced.

# Row 2
$ws.Range("D2").Value = '88.490.20'
$ws.Range("E2").Value = '  -2.61%  '

# Row 3
$ws.Range("D3").Value = '3.122.18'
$ws.Range("E3").Value = '  -1.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '634.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.391'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.79%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.777'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +13.15%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.05%  '

# Row 10
$ws.Range("D10").Value = '3.121.31'
$ws.Range("E10").Value = '  -1.79%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.567'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.40%  '

# Row 12
$ws.Range("E12").Value = '  +1.36%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.83%  '

# Row 15
$ws.Range("D15").Value = '88.344.58'
$ws.Range("E15").Value = '  -2.36%  '

# Row 16
$ws.Range("D16").Value = '3.693.41'
$ws.Range("E16").Value = '  -1.71%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.49%  '

# Row 18
$ws.Range("D18").Value = '3.129.45'
$ws.Range("E18").Value = '  -1.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.00%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000222'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +17.00%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.69%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '422.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.80%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.68%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.10%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.15%  '

# Row 28
$ws.Range("D28").Value = '3.283.52'
$ws.Range("E28").Value = '  -2.17%  '

# Row 29
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.03%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.157'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.33%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.14%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.92'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.56%  '

# Row 37
$ws.Range("E37").Value = '  +2.26%  '

# Row 38
$ws.Range("E38").Value = '  -1.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.98%  '

# Row 40
$ws.Range("E40").Value = '  -0.45%  '

# Row 41
$ws.Range("E41").Value = '  +0.35%  '

# Row 42
$ws.Range("E42").Value = '  +0.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.366'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.35%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.82%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.76%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.132'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.76%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.87%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0654'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.49%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '162.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.716'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '

# Row 34: Bittensor -> Kaspa (coin identity changes, new price/volume)
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.148"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.71%  "

# Row 35: Kaspa -> Bittensor (coin identity changes, new price/volume)
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "502.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.33%  "

